# Update PutShipCommand sequence diagram.
# Update putShip() command call in sequence diagram.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Helper: find a top-level shape on the slide by its PowerPoint shape Id.
function Get-ShapeById($slide, $id) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $sh = $slide.Shapes.Item($i)
        if ($sh.Id -eq $id) {
            return $sh
        }
    }
    return $null
}

$EMU_PER_PT = 12700

# ---------------------------------------------------------------------
# 1. Resize/reposition the two big background "lane" rectangles.
# ---------------------------------------------------------------------

# id=170 "Rectangle 65" (right-hand "Model" lane background)
$sh170 = Get-ShapeById $s 170
$sh170.Left = 9699812 / $EMU_PER_PT
$sh170.Top = 296676 / $EMU_PER_PT
$sh170.Width = 7075272 / $EMU_PER_PT
$sh170.Height = 10908136 / $EMU_PER_PT

# id=171 "Rectangle 65" (big outer lane background)
$sh171 = Get-ShapeById $s 171
$sh171.Width = 9336305 / $EMU_PER_PT
$sh171.Height = 10908137 / $EMU_PER_PT

# ---------------------------------------------------------------------
# 2. Extend the lifeline connectors so they reach the (now shorter) bottom.
# ---------------------------------------------------------------------

# id=46 "Straight Connector 45" (:Model lifeline)
$sh46 = Get-ShapeById $s 46
$sh46.Height = 7985954 / $EMU_PER_PT

# id=188 "Straight Connector 187" (putShip actor lifeline)
$sh188 = Get-ShapeById $s 188
$sh188.Height = 6043459 / $EMU_PER_PT

# id=110 "Straight Connector 109" (:Player lifeline)
$sh110 = Get-ShapeById $s 110
$sh110.Height = 7780974 / $EMU_PER_PT

# id=126 "Straight Connector 125" (:Checker lifeline)
$sh126 = Get-ShapeById $s 126
$sh126.Height = 6848653 / $EMU_PER_PT

# ---------------------------------------------------------------------
# 3. Move the "args1/args2" caption text box up to follow the new bottom.
# ---------------------------------------------------------------------

# id=61 "TextBox 60"
$sh61 = Get-ShapeById $s 61
$sh61.Top = 11361872 / $EMU_PER_PT

# ---------------------------------------------------------------------
# 4. Ungroup "Group 134" (id=135): it contained the putShip() call-out
#    rectangle + label plus two small connector glyphs. Drop the small
#    glyphs, keep the rectangle + label as independent top-level shapes,
#    and move the label to its new spot near the Model lane.
# ---------------------------------------------------------------------

$grp135 = Get-ShapeById $s 135
$ungrouped = $grp135.Ungroup()

# After ungrouping, delete the now-unneeded "Rectangle 63" (id=64) and
# the two curved connector glyphs (id=74, id=88).
$idsToRemove = @(64, 74, 88)
foreach ($id in $idsToRemove) {
    $victim = Get-ShapeById $s $id
    if ($victim -ne $null) {
        $victim.Delete()
    }
}

# Reposition the surviving label textbox (id=68 "TextBox 67", "putShip()")
# to its new location.
$sh68 = Get-ShapeById $s 68
$sh68.Left = 9185906 / $EMU_PER_PT
$sh68.Top = 6872592 / $EMU_PER_PT
